$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-15: update B (name), C, D, E values to reflect the new line7/line8
# insertion and re-balanced contingency results.
$rows = @(
    @{ Row = 8;  Name = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  Name = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; Name = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; Name = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; Name = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; Name = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; Name = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; Name = "extr6"; C = 7;  D = 11; E = $false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}

# New rows 16-17: extr7, extr8 (shifted down from the old extr7/extr8 rows
# which used to be rows 14-15).
$newRows = @(
    @{ Row = 16; Id = 14; Name = "extr7"; C = 5; D = 7; E = $false },
    @{ Row = 17; Id = 15; Name = "extr8"; C = 8; D = 5; E = $false }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    # Match the bold / thin-border / centered style used by the rest of
    # column A (the "name index" column).
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}
